# Auto-generated edit script
# Applies numeric cell updates per the Fenrir_Profits.xlsx commit diff
# (scheduled market-price refresh: H/I/J/K/L = price stats, M/N = profit margins)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 673.86664
$ws.Range("I18").Value = 685.2308
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 685.2308
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -401.2308
$ws.Range("N18").Value = -1168
# Row 58
$ws.Range("H58").Value = 728
$ws.Range("I58").Value = 365.8889
$ws.Range("J58").Value = 931.6875
$ws.Range("K58").Value = 1097.6667
$ws.Range("L58").Value = 2795.0625
$ws.Range("M58").Value = -947.6667
$ws.Range("N58").Value = -3095.0625
# Row 74
$ws.Range("H74").Value = 5461.222
$ws.Range("I74").Value = 6164.4287
$ws.Range("K74").Value = 6164.4287
$ws.Range("M74").Value = -5228.4287
# Row 76
$ws.Range("H76").Value = 4389.3228
$ws.Range("I76").Value = 3817.2856
$ws.Range("K76").Value = 3817.2856
$ws.Range("M76").Value = -3502.2856
# Row 77
$ws.Range("H77").Value = 5461.222
$ws.Range("I77").Value = 6164.4287
$ws.Range("K77").Value = 30822.1435
$ws.Range("M77").Value = -26142.1435
# Row 79
$ws.Range("H79").Value = 4389.3228
$ws.Range("I79").Value = 3817.2856
$ws.Range("K79").Value = 3817.2856
$ws.Range("M79").Value = -2725.2856
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 129
$ws.Range("H129").Value = 49548.125
$ws.Range("I129").Value = 277.7143
$ws.Range("J129").Value = 87869.55499999999
$ws.Range("K129").Value = 833.1428999999999
$ws.Range("L129").Value = 263608.665
$ws.Range("M129").Value = 4166.8571
$ws.Range("N129").Value = -273608.665
# Row 132
$ws.Range("H132").Value = 42087830
$ws.Range("I132").Value = 63129224
$ws.Range("K132").Value = 189387672
$ws.Range("M132").Value = -189385142
# Row 135
$ws.Range("H135").Value = 2187.8572
$ws.Range("I135").Value = 680.69446
$ws.Range("J135").Value = 6361.5386
$ws.Range("K135").Value = 6126.25014
$ws.Range("L135").Value = 57253.8474
$ws.Range("M135").Value = -3591.25014
$ws.Range("N135").Value = -62323.8474
# Row 137
$ws.Range("H137").Value = 1147.8644
$ws.Range("I137").Value = 963.25
$ws.Range("J137").Value = 1536.5264
$ws.Range("K137").Value = 2889.75
$ws.Range("L137").Value = 4609.5792
$ws.Range("M137").Value = -339.75
$ws.Range("N137").Value = -9709.5792
# Row 138
$ws.Range("H138").Value = 3625.2805
$ws.Range("I138").Value = 4710.8096
$ws.Range("J138").Value = 2485.475
$ws.Range("K138").Value = 14132.4288
$ws.Range("L138").Value = 7456.424999999999
$ws.Range("M138").Value = -8992.428799999998
$ws.Range("N138").Value = -17736.425

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4124.4736
$ws.Range("I61").Value = 4150.0586
$ws.Range("J61").Value = 3907
$ws.Range("K61").Value = 4150.0586
$ws.Range("L61").Value = 3907
$ws.Range("M61").Value = -3938.0586
$ws.Range("N61").Value = -4331
# Row 63
$ws.Range("H63").Value = 1854391.4
$ws.Range("I63").Value = 3705522.5
$ws.Range("J63").Value = 3260.3333
$ws.Range("K63").Value = 3705522.5
$ws.Range("L63").Value = 3260.3333
$ws.Range("M63").Value = -3704836.5
$ws.Range("N63").Value = -4632.3333
# Row 66
$ws.Range("H66").Value = 1854391.4
$ws.Range("I66").Value = 3705522.5
$ws.Range("J66").Value = 3260.3333
$ws.Range("K66").Value = 18527612.5
$ws.Range("L66").Value = 16301.6665
$ws.Range("M66").Value = -18524180.5
$ws.Range("N66").Value = -23165.6665
# Row 74
$ws.Range("H74").Value = 911.8570999999999
$ws.Range("I74").Value = 722.9091
$ws.Range("K74").Value = 722.9091
$ws.Range("M74").Value = 151.0909
# Row 77
$ws.Range("H77").Value = 911.8570999999999
$ws.Range("I77").Value = 722.9091
$ws.Range("K77").Value = 3614.5455
$ws.Range("M77").Value = 753.4545000000003
# Row 88
$ws.Range("H88").Value = 15137.637
$ws.Range("I88").Value = 3400
$ws.Range("J88").Value = 17746
$ws.Range("K88").Value = 3400
$ws.Range("L88").Value = 17746
$ws.Range("M88").Value = -2994
$ws.Range("N88").Value = -18558
# Row 91
$ws.Range("H91").Value = 15137.637
$ws.Range("I91").Value = 3400
$ws.Range("J91").Value = 17746
$ws.Range("K91").Value = 3400
$ws.Range("L91").Value = 17746
$ws.Range("M91").Value = -1996
$ws.Range("N91").Value = -20554
# Row 136
$ws.Range("H136").Value = 4124.4736
$ws.Range("I136").Value = 4150.0586
$ws.Range("J136").Value = 3907
$ws.Range("K136").Value = 12450.1758
$ws.Range("L136").Value = 11721
$ws.Range("M136").Value = -9900.175800000001
$ws.Range("N136").Value = -16821

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 8335035.5
$ws.Range("I86").Value = 25001500
$ws.Range("K86").Value = 25001500
$ws.Range("M86").Value = -25000377
# Row 89
$ws.Range("H89").Value = 8335035.5
$ws.Range("I89").Value = 25001500
$ws.Range("K89").Value = 125007500
$ws.Range("M89").Value = -125001884

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7521.355
$ws.Range("I31").Value = 2736.3965
$ws.Range("J31").Value = 76903.25
$ws.Range("K31").Value = 2736.3965
$ws.Range("L31").Value = 76903.25
$ws.Range("M31").Value = -2441.3965
$ws.Range("N31").Value = -77493.25
# Row 34
$ws.Range("H34").Value = 7521.355
$ws.Range("I34").Value = 2736.3965
$ws.Range("J34").Value = 76903.25
$ws.Range("K34").Value = 2736.3965
$ws.Range("L34").Value = 76903.25
$ws.Range("M34").Value = -2534.3965
$ws.Range("N34").Value = -77307.25
# Row 132
$ws.Range("H132").Value = 8776428
$ws.Range("I132").Value = 16667620
$ws.Range("J132").Value = 8437.833000000001
$ws.Range("K132").Value = 50002860
$ws.Range("L132").Value = 25313.499
$ws.Range("M132").Value = -50000330
$ws.Range("N132").Value = -30373.499
# Row 134
$ws.Range("H134").Value = 6650016
$ws.Range("I134").Value = 6411163
$ws.Range("J134").Value = 7814425
$ws.Range("K134").Value = 19233489
$ws.Range("L134").Value = 23443275
$ws.Range("M134").Value = -19230954
$ws.Range("N134").Value = -23448345

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 150.15152
$ws.Range("I12").Value = 335.77777
$ws.Range("J12").Value = 80.541664
$ws.Range("K12").Value = 1007.33331
$ws.Range("L12").Value = 241.624992
$ws.Range("M12").Value = -834.33331
$ws.Range("N12").Value = -587.624992

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 34451.21
$ws.Range("I70").Value = 51175.715
$ws.Range("K70").Value = 51175.715
$ws.Range("M70").Value = -50905.715
# Row 73
$ws.Range("H73").Value = 34451.21
$ws.Range("I73").Value = 51175.715
$ws.Range("K73").Value = 51175.715
$ws.Range("M73").Value = -50239.715

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 5298.125
$ws.Range("I82").Value = 2999.5
$ws.Range("J82").Value = 6064.3335
$ws.Range("K82").Value = 2999.5
$ws.Range("L82").Value = 6064.3335
$ws.Range("M82").Value = -2638.5
$ws.Range("N82").Value = -6786.3335
# Row 85
$ws.Range("H85").Value = 5298.125
$ws.Range("I85").Value = 2999.5
$ws.Range("J85").Value = 6064.3335
$ws.Range("K85").Value = 2999.5
$ws.Range("L85").Value = 6064.3335
$ws.Range("M85").Value = -1751.5
$ws.Range("N85").Value = -8560.333500000001

